$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stale "_GoBack" bookmark that currently sits in the
#    empty paragraph right after the dependency-injection sentence.
#    (it gets re-created at a new location below)
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) "3 Decorator pattern" -> split the leading "3 " run into "3"
#    plus a separate " " run, with a (new) "_GoBack" bookmark sitting
#    between them.
# ------------------------------------------------------------------
$rDecorator = $d.Content
$rDecorator.Find.Execute("Decorator pattern")
$decoratorStart = $rDecorator.Start

# force a run boundary right before the trailing space of "3 " by
# nudging a formatting property and then reverting it
$rSpace1 = $d.Range($decoratorStart - 1, $decoratorStart)
$rSpace1.Font.Size = 12.5
$rSpace1.Font.Size = 12

# drop the (new) bookmark exactly between "3" and the space
$bmRange = $d.Range($decoratorStart - 1, $decoratorStart - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 3) "Facade pattern" -> "Repository pattern", ending up as three
#    separate runs: "Repository", " ", "pattern".
# ------------------------------------------------------------------
$rFacadeLine = $d.Content
$rFacadeLine.Find.Execute("Facade pattern")
$facadeStart = $rFacadeLine.Start

$rFacadeWord = $d.Range($facadeStart, $facadeStart + 6)
$rFacadeWord.Text = "Repository"

$rMidSpace = $d.Range($facadeStart + 10, $facadeStart + 11)
$rMidSpace.Font.Size = 12.5
$rMidSpace.Font.Size = 12

$rPatternWord = $d.Range($facadeStart + 11, $facadeStart + 18)
$rPatternWord.Font.Size = 12.5
$rPatternWord.Font.Size = 12
